$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Input")
$ws2 = $wb.Worksheets.Item("Output")

# --- Cell value updates (order chosen so new shared-string entries come
# out in the same sequence as the target file) ---
$ws1.Range("A8").Value = "GroupAddClient"
$ws1.Range("B7").Value = "Jhon Deer"
$ws1.Range("B8").Value = "click"

$ws1.Range("B1").Copy()
$ws2.Range("B1").PasteSpecial(-4122)

$ws2.Range("B1").Value = "Group3"
$ws2.Range("A2").Value = "VerifyClientCreated"
$ws2.Range("B2").Value = "Jhon Deer"

# --- Selection / active-sheet updates ---
[void]$ws1.Range("B2").Select()
[void]$ws2.Activate()
[void]$ws2.Range("C15").Select()

# --- Column width update on Output sheet ---
$ws2.Range("A1").ColumnWidth = 16.166666666666668
